# Update '想去人数' (F column) numeric values per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 199
$ws.Range("F5").Value = 977
$ws.Range("F6").Value = 5404
$ws.Range("F7").Value = 470
$ws.Range("F9").Value = 935
$ws.Range("F11").Value = 75
$ws.Range("F12").Value = 33
$ws.Range("F13").Value = 580
$ws.Range("F17").Value = 1801
$ws.Range("F19").Value = 869
$ws.Range("F22").Value = 322
$ws.Range("F24").Value = 142
$ws.Range("F25").Value = 1052
$ws.Range("F28").Value = 2752
$ws.Range("F31").Value = 63
$ws.Range("F32").Value = 108
$ws.Range("F33").Value = 31
$ws.Range("F34").Value = 345
$ws.Range("F39").Value = 283
$ws.Range("F40").Value = 673
$ws.Range("F41").Value = 84
$ws.Range("F44").Value = 64

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 172
$ws.Range("F6").Value = 120

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 236

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 236
$ws.Range("F3").Value = 199
$ws.Range("F5").Value = 977
$ws.Range("F7").Value = 5404
$ws.Range("F8").Value = 470
$ws.Range("F11").Value = 172
$ws.Range("F12").Value = 935
$ws.Range("F15").Value = 120
$ws.Range("F16").Value = 75
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 580
$ws.Range("F23").Value = 1801
$ws.Range("F25").Value = 869
$ws.Range("F27").Value = 322
$ws.Range("F30").Value = 142
$ws.Range("F31").Value = 1052
$ws.Range("F32").Value = 2752
$ws.Range("F35").Value = 63
$ws.Range("F36").Value = 108
$ws.Range("F37").Value = 31
$ws.Range("F38").Value = 345
$ws.Range("F42").Value = 283
$ws.Range("F43").Value = 673
$ws.Range("F44").Value = 84
$ws.Range("F46").Value = 64
